$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: rows 57-60 - the match data (columns F:V) got cyclically
# reshuffled: the old row 58 data now lives in row 57, old row 59 data
# moved into row 58, old row 60 data moved into row 59, and the old row 57
# data moved down into row 60. Columns A:E (index/country/tournament/
# season/date) are unchanged.
$row57 = $ws.Range("F57:V57").Value2
$row58 = $ws.Range("F58:V58").Value2
$row59 = $ws.Range("F59:V59").Value2
$row60 = $ws.Range("F60:V60").Value2

$ws.Range("F57:V57").Value = $row58
$ws.Range("F58:V58").Value = $row59
$ws.Range("F59:V59").Value = $row60
$ws.Range("F60:V60").Value = $row57

# --- Part 2: append a new match record as row 83 ---
# Duplicate formatting/styles from the last existing row (82) into the new
# row 83, then overwrite with the new match's values.
$ws.Range("A82:V82").Copy($ws.Range("A83:V83"))

$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "germany"
$ws.Range("C83").Value = "bundesliga"
$ws.Range("D83").Value = "2023-2024"
$ws.Range("E83").Value = 45233.85416666666
$ws.Range("F83").Value = "Darmstadt"
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = "Bochum"
$ws.Range("I83").Value = 2
$ws.Range("J83").Value = 2.44
$ws.Range("K83").Value = "22/10/2023 15:56"
$ws.Range("L83").Value = 2.37
$ws.Range("M83").Value = "03/11/2023 19:57"
$ws.Range("N83").Value = 3.55
$ws.Range("O83").Value = "22/10/2023 15:56"
$ws.Range("P83").Value = 3.53
$ws.Range("Q83").Value = "03/11/2023 19:44"
$ws.Range("R83").Value = 2.92
$ws.Range("S83").Value = "22/10/2023 15:56"
$ws.Range("T83").Value = 3.1
$ws.Range("U83").Value = "03/11/2023 20:00"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/germany/bundesliga/darmstadt-bochum/82olJ6dE/"
